$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Refresh the cached "datetimeFigureOut" date field text (6/27/2017 ->
#    6/30/2017) on the slide master and every slide layout's Date
#    Placeholder shape.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$newDate = "6/30/2017"

$masterShapes = $master.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $shp = $masterShapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $lay = $master.CustomLayouts.Item($li)
    $layShapes = $lay.Shapes
    for ($i = 1; $i -le $layShapes.Count; $i++) {
        $shp = $layShapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Re-color the Timer shapes on slide 1 to the shared palette:
#      accent5 lumMod 75% -> 95B3D7   (TimerBody fill + line)
#      accent1 lumMod 60%/lumOff 40% -> 558ED5 (TimerLineMarker1..5 line)
#      accent5 lumMod 50% -> 0070C0   (TimerSliderHead / TimerSliderBody fill)
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)
$shapes = $s.Shapes

for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)

    if ($shp.Name -eq "TimerBody") {
        $shp.Fill.ForeColor.RGB = 14136213
        $shp.Line.ForeColor.RGB = 14136213
    }

    if ($shp.Name -eq "TimerSliderHead") {
        $shp.Fill.ForeColor.RGB = 12611584
    }

    if ($shp.Name -eq "TimerSliderBody") {
        $shp.Fill.ForeColor.RGB = 12611584
    }

    if ($shp.Name -eq "TimerLineMarkerGroup") {
        $grpItems = $shp.GroupItems
        for ($j = 1; $j -le $grpItems.Count; $j++) {
            $sub = $grpItems.Item($j)
            if ($sub.Name -like "TimerLineMarker*") {
                $sub.Line.ForeColor.RGB = 13995605
            }
        }
    }
}
